$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is numeric-looking (e.g. "7.49") need the column
# format forced to Text first, otherwise Excel would store them as real numbers
# (losing the original text formatting, e.g. trailing zeros like "18.90").
$textForced = @("D5", "D6", "D10", "D11", "D12", "D13", "D17", "D19", "D21", "D23", "D24", "D26", "D28", "D29", "D30", "D32", "D33", "D38", "D39", "D44", "D45", "D50", "D51")
foreach ($addr in $textForced) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '51.848.08'
$ws.Range("E2").Value = '  -0.55%  '
$ws.Range("D3").Value = '2.781.49'
$ws.Range("E3").Value = '  -2.04%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '356.53'
$ws.Range("E5").Value = '  -1.65%  '
$ws.Range("D6").Value = '109.13'
$ws.Range("E6").Value = '  -3.47%  '
$ws.Range("E7").Value = '  -2.64%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -2.31%  '
$ws.Range("D10").Value = '40.18'
$ws.Range("E10").Value = '  -3.52%  '
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").Value = '  -1.89%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.134'
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").Value = '19.42'
$ws.Range("E13").Value = '  -3.36%  '
$ws.Range("E14").Value = '  -3.22%  '
$ws.Range("D15").Value = '3.225.75'
$ws.Range("E15").Value = '  -1.81%  '
$ws.Range("D16").Value = '2.766.24'
$ws.Range("E16").Value = '  -2.40%  '
$ws.Range("D17").Value = '0.937'
$ws.Range("E17").Value = '  +2.65%  '
$ws.Range("D18").Value = '51.820.87'
$ws.Range("E18").Value = '  -0.49%  '
$ws.Range("D19").Value = '7.49'
$ws.Range("E19").Value = '  -0.97%  '
$ws.Range("E20").Value = '  -2.44%  '
$ws.Range("D21").Value = '13.09'
$ws.Range("E21").Value = '  -3.62%  '
$ws.Range("E22").Value = '  -2.53%  '
$ws.Range("D23").Value = '70.03'
$ws.Range("E23").Value = '  -0.48%  '
$ws.Range("D24").Value = '269.53'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("E25").Value = '  -3.63%  '
$ws.Range("D26").Value = '26.46'
$ws.Range("E26").Value = '  -2.36%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("D28").Value = '0.163'
$ws.Range("E28").Value = '  +16.35%  '
$ws.Range("D29").Value = '10.29'
$ws.Range("E29").Value = '  -1.00%  '
$ws.Range("D30").Value = '2.13'
$ws.Range("E30").Value = '  -5.25%  '
$ws.Range("E31").Value = '  -3.87%  '
$ws.Range("D32").Value = '51.99'
$ws.Range("E32").Value = '  -3.37%  '
$ws.Range("D33").Value = '34.36'
$ws.Range("E33").Value = '  -3.39%  '
$ws.Range("E34").Value = '  -2.51%  '
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("E36").Value = '  -5.03%  '
$ws.Range("E37").Value = '  +0.00%  '
$ws.Range("D38").Value = '18.90'
$ws.Range("E38").Value = '  +2.38%  '
$ws.Range("D39").Value = '3.19'
$ws.Range("E39").Value = '  -2.80%  '
$ws.Range("E40").Value = '  -4.47%  '
$ws.Range("E41").Value = '  +3.17%  '
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("D44").Value = '119.54'
$ws.Range("E44").Value = '  -5.78%  '
$ws.Range("D45").Value = '21.81'
$ws.Range("E45").Value = '  -7.71%  '
$ws.Range("D46").Value = '2.087.63'
$ws.Range("E46").Value = '  -1.27%  '
$ws.Range("E47").Value = '  -4.89%  '
$ws.Range("E48").Value = '  -1.88%  '
$ws.Range("E49").Value = '  -2.57%  '
$ws.Range("D50").Value = '0.956'
$ws.Range("E50").Value = '  -3.28%  '
$ws.Range("B51").Value = 'BitgetToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$ws.Range("D51").Value = '1.13'
$ws.Range("E51").Value = '  +31.10%  '
